# Updates cryptos list prices/volumes and swaps TheGraph/PEPE rows (38 and 39)
# Generated from the canonical OOXML diff for the Mon Mar 11 19:11:33 UTC 2024 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Addr = "D2"; Value = "72.688.82" },
    @{ Addr = "E2"; Value = "  +4.97%  " },
    @{ Addr = "D3"; Value = "4.057.70" },
    @{ Addr = "E3"; Value = "  +4.14%  " },
    @{ Addr = "E4"; Value = "  +0.12%  " },
    @{ Addr = "D5"; Value = "'519.60" },
    @{ Addr = "E5"; Value = "  -1.82%  " },
    @{ Addr = "D6"; Value = "'147.06" },
    @{ Addr = "E6"; Value = "  +1.67%  " },
    @{ Addr = "D7"; Value = "'0.740" },
    @{ Addr = "E7"; Value = "  +20.88%  " },
    @{ Addr = "D8"; Value = "4.045.07" },
    @{ Addr = "E8"; Value = "  +4.31%  " },
    @{ Addr = "E9"; Value = "  +0.09%  " },
    @{ Addr = "D10"; Value = "'0.774" },
    @{ Addr = "E10"; Value = "  +7.99%  " },
    @{ Addr = "D11"; Value = "'0.177" },
    @{ Addr = "E11"; Value = "  +2.47%  " },
    @{ Addr = "D12"; Value = "'0.0000328" },
    @{ Addr = "E12"; Value = "  -2.47%  " },
    @{ Addr = "D13"; Value = "'47.76" },
    @{ Addr = "E13"; Value = "  +13.54%  " },
    @{ Addr = "D14"; Value = "'11.16" },
    @{ Addr = "E14"; Value = "  +9.01%  " },
    @{ Addr = "D15"; Value = "4.708.32" },
    @{ Addr = "E15"; Value = "  +4.11%  " },
    @{ Addr = "D16"; Value = "4.066.72" },
    @{ Addr = "E16"; Value = "  +3.83%  " },
    @{ Addr = "D17"; Value = "'21.28" },
    @{ Addr = "E17"; Value = "  +7.64%  " },
    @{ Addr = "D18"; Value = "'14.15" },
    @{ Addr = "E18"; Value = "  +0.54%  " },
    @{ Addr = "D19"; Value = "'1.21" },
    @{ Addr = "E19"; Value = "  -0.06%  " },
    @{ Addr = "D21"; Value = "72.500.74" },
    @{ Addr = "E21"; Value = "  +4.77%  " },
    @{ Addr = "D22"; Value = "'442.21" },
    @{ Addr = "E22"; Value = "  +3.97%  " },
    @{ Addr = "D23"; Value = "'104.76" },
    @{ Addr = "E23"; Value = "  +18.73%  " },
    @{ Addr = "D24"; Value = "'3.58" },
    @{ Addr = "E24"; Value = "  +6.02%  " },
    @{ Addr = "D25"; Value = "'14.81" },
    @{ Addr = "E25"; Value = "  +4.87%  " },
    @{ Addr = "D26"; Value = "'4.01" },
    @{ Addr = "E26"; Value = "  -0.49%  " },
    @{ Addr = "D27"; Value = "'11.45" },
    @{ Addr = "E27"; Value = "  +0.96%  " },
    @{ Addr = "D28"; Value = "'11.03" },
    @{ Addr = "E28"; Value = "  +4.03%  " },
    @{ Addr = "D29"; Value = "'37.75" },
    @{ Addr = "E29"; Value = "  +3.72%  " },
    @{ Addr = "E30"; Value = "  +2.49%  " },
    @{ Addr = "D31"; Value = "'3.29" },
    @{ Addr = "E31"; Value = "  +16.68%  " },
    @{ Addr = "D32"; Value = "'13.70" },
    @{ Addr = "E32"; Value = "  +4.47%  " },
    @{ Addr = "D33"; Value = "'0.130" },
    @{ Addr = "E33"; Value = "  +3.84%  " },
    @{ Addr = "D34"; Value = "'677.63" },
    @{ Addr = "E34"; Value = "  -0.46%  " },
    @{ Addr = "D35"; Value = "'6.83" },
    @{ Addr = "E35"; Value = "  +15.07%  " },
    @{ Addr = "D36"; Value = "'67.19" },
    @{ Addr = "E36"; Value = "  -2.18%  " },
    @{ Addr = "D37"; Value = "'42.83" },
    @{ Addr = "E37"; Value = "  +6.89%  " },
    @{ Addr = "B38"; Value = "PEPE" },
    @{ Addr = "C38"; Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe" },
    @{ Addr = "D38"; Value = "0.0₃0866" },
    @{ Addr = "E38"; Value = "  -0.51%  " },
    @{ Addr = "B39"; Value = "TheGraph" },
    @{ Addr = "C39"; Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt" },
    @{ Addr = "D39"; Value = "'0.429" },
    @{ Addr = "E39"; Value = "  -1.07%  " },
    @{ Addr = "D40"; Value = "'3.54" },
    @{ Addr = "E40"; Value = "  +8.54%  " },
    @{ Addr = "D41"; Value = "'0.151" },
    @{ Addr = "E41"; Value = "  +1.63%  " },
    @{ Addr = "E42"; Value = "  +0.12%  " },
    @{ Addr = "D43"; Value = "'0.0498" },
    @{ Addr = "E43"; Value = "  +3.71%  " },
    @{ Addr = "E44"; Value = "  -0.24%  " },
    @{ Addr = "D45"; Value = "'3.27" },
    @{ Addr = "E45"; Value = "  +2.85%  " },
    @{ Addr = "D46"; Value = "'0.158" },
    @{ Addr = "E46"; Value = "  +12.88%  " },
    @{ Addr = "D47"; Value = "'2.70" },
    @{ Addr = "E47"; Value = "  -2.98%  " },
    @{ Addr = "D48"; Value = "'3.51" },
    @{ Addr = "E48"; Value = "  +4.01%  " },
    @{ Addr = "E49"; Value = "  +2.55%  " },
    @{ Addr = "D50"; Value = "'9.14" },
    @{ Addr = "E50"; Value = "  +7.30%  " },
    @{ Addr = "D51"; Value = "'3.33" },
    @{ Addr = "E51"; Value = "  +2.39%  " }
)

foreach ($u in $updates) {
    $ws.Range($u.Addr).Value = $u.Value
}
